$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "Sin precio"
$ws.Range("E2").Value = ""
$ws.Range("F3").Value = "`$ 1.758,58"
$ws.Range("F4").Value = "`$ 2.036,27"
$ws.Range("F5").Value = "`$ 1.527,18"
$ws.Range("F6").Value = "`$ 629,33"
$ws.Range("F8").Value = "`$ 1.758,58"
$ws.Range("F9").Value = "`$ 277,59"
$ws.Range("F10").Value = "`$ 277,59"
$ws.Range("F12").Value = "`$ 509,00"
$ws.Range("F13").Value = "`$ 1.110,65"
$ws.Range("F14").Value = "`$ 2.545,36"
$ws.Range("F15").Value = "`$ 2.545,36"
$ws.Range("F16").Value = "`$ 2.545,36"
$ws.Range("F17").Value = "`$ 3.424,69"
$ws.Range("F18").Value = "`$ 1.249,49"
$ws.Range("F19").Value = "`$ 1.851,15"
$ws.Range("F20").Value = "`$ 1.110,65"
$ws.Range("F21").Value = "`$ 925,52"
$ws.Range("F22").Value = "`$ 1.851,15"
$ws.Range("F24").Value = "`$ 3.054,45"
$ws.Range("F25").Value = "`$ 1.527,18"
$ws.Range("F26").Value = "`$ 2.684,20"
$ws.Range("F27").Value = "`$ 370,16"
$ws.Range("F28").Value = "`$ 416,44"
$ws.Range("F29").Value = "`$ 416,44"
$ws.Range("F30").Value = "`$ 370,16"
$ws.Range("F32").Value = "`$ 2.313,95"
$ws.Range("F33").Value = "`$ 2.313,95"
$ws.Range("F34").Value = "`$ 1.418,89"
$ws.Range("F35").Value = "`$ 1.216,19"
$ws.Range("F36").Value = "`$ 1.368,23"
$ws.Range("F37").Value = "`$ 462,72"
$ws.Range("F38").Value = "`$ 1.203,22"
$ws.Range("F39").Value = "`$ 490,48"
$ws.Range("F40").Value = "`$ 1.156,94"
$ws.Range("F41").Value = "`$ 1.240,23"
$ws.Range("F42").Value = "`$ 444,20"
$ws.Range("F43").Value = "`$ 509,00"
$ws.Range("F44").Value = "`$ 925,52"
$ws.Range("F45").Value = "Sin precio"
$ws.Range("E45").Value = ""
$ws.Range("F46").Value = "`$ 3.239,58"
$ws.Range("F47").Value = "`$ 1.388,34"
$ws.Range("F48").Value = "`$ 573,79"
$ws.Range("F49").Value = "`$ 2.776,76"
$ws.Range("F50").Value = "`$ 1.758,58"
$ws.Range("F51").Value = "`$ 573,79"
$ws.Range("F52").Value = "`$ 768,17"
$ws.Range("F53").Value = "`$ 462,72"
$ws.Range("F54").Value = "`$ 879,24"
$ws.Range("F55").Value = "`$ 2.036,26"
$ws.Range("F56").Value = "`$ 569,16"
$ws.Range("F57").Value = "`$ 569,16"
$ws.Range("F58").Value = "`$ 1.110,65"
$ws.Range("F59").Value = "Sin precio"
$ws.Range("E59").Value = ""
$ws.Range("F60").Value = "Sin precio"
$ws.Range("E60").Value = ""
$ws.Range("F61").Value = "`$ 462,72"
$ws.Range("F62").Value = "Sin precio"
$ws.Range("E62").Value = ""
$ws.Range("F65").Value = "`$ 490,48"
$ws.Range("F66").Value = "`$ 490,48"
$ws.Range("F67").Value = "`$ 1.511,89"
$ws.Range("F68").Value = "`$ 1.388,69"
$ws.Range("F69").Value = "`$ 1.119,89"
$ws.Range("F70").Value = "`$ 1.455,89"
$ws.Range("F71").Value = "`$ 1.175,89"
$ws.Range("F72").Value = "`$ 1.036,59"
$ws.Range("F73").Value = "`$ 647,84"
$ws.Range("F74").Value = "`$ 2.036,27"
$ws.Range("F75").Value = "`$ 2.406,51"
$ws.Range("F77").Value = "`$ 1.156,93"
$ws.Range("F78").Value = "`$ 1.156,93"
$ws.Range("F79").Value = "`$ 1.156,93"
$ws.Range("F80").Value = "`$ 1.156,93"
$ws.Range("F81").Value = "`$ 1.156,93"
$ws.Range("F82").Value = "`$ 2.684,20"
$ws.Range("F83").Value = "`$ 2.684,20"
$ws.Range("F84").Value = "`$ 2.684,20"
$ws.Range("F85").Value = "`$ 2.684,20"
$ws.Range("F86").Value = "`$ 1.110,65"
$ws.Range("F87").Value = "`$ 2.684,20"
